$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (currency / percentage) cells that must remain text, matching existing formatting
$textChanges = @{
    "B3" = "$-3963.69"
    "C3" = "$-1614.28"
    "D3" = "$-1442.63"
    "E3" = "$-1267.01"
    "F3" = "$-1087.32"
    "G3" = "$-903.48"
    "H3" = "$-715.39"
    "I3" = "$-522.96"
    "J3" = "$-326.11"
    "K3" = "$-124.72"
    "B6" = "$21764.66"
    "C6" = "$23976.13"
    "D6" = "$26497.84"
    "E6" = "$29341.11"
    "F6" = "$32517.83"
    "G6" = "$36040.51"
    "H6" = "$39922.29"
    "I6" = "$44176.99"
    "J6" = "$48819.18"
    "K6" = "$53864.20"
    "B7" = "-9.04%"
    "C7" = "-3.68%"
    "D7" = "-3.29%"
    "E7" = "-2.89%"
    "F7" = "-2.48%"
    "G7" = "-2.06%"
    "H7" = "-1.63%"
    "I7" = "-1.19%"
    "J7" = "-0.74%"
    "K7" = "-0.28%"
    "C8" = "-45.34%"
    "D8" = "-22.28%"
    "E8" = "-12.54%"
    "F8" = "-7.21%"
    "G8" = "-3.85%"
    "H8" = "-1.56%"
    "I8" = "0.10%"
    "J8" = "1.35%"
    "K8" = "2.31%"
    "B9" = "3.33%"
    "C9" = "4.95%"
    "D9" = "5.07%"
    "E9" = "5.19%"
    "F9" = "5.31%"
    "G9" = "5.44%"
    "H9" = "5.57%"
    "I9" = "5.70%"
    "J9" = "5.84%"
    "K9" = "5.98%"
}

foreach ($cell in $textChanges.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $textChanges[$cell]
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Numeric cells (row 11, Expense/rent ratio)
$numChanges = @{
    "C11" = 1.1
    "D11" = 1.09
    "E11" = 1.08
    "F11" = 1.06
    "G11" = 1.05
    "H11" = 1.04
    "I11" = 1.03
    "J11" = 1.02
    "K11" = 1.01
}

foreach ($cell in $numChanges.Keys) {
    $ws.Range($cell).Value = $numChanges[$cell]
}
